$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf16"
$ws.Range("C2").Value = "Fgfr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1628603333333333
$ws.Range("H2").Value = 0.488581
$ws.Range("I2").Value = 0.06904471801498467
$ws.Range("J2").Value = 0.06904471801498467
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.656156333333334
$ws.Range("N2").Value = 4.968469000000001
$ws.Range("O2").Value = 0.6151212440816572
$ws.Range("P2").Value = 0.6151212440816572
$ws.Range("Q2").Value = 0.2697221724987778
$ws.Range("R2").Value = 2.427499552489
$ws.Range("S2").Value = 0.04247087284264459
$ws.Range("T2").Value = 0.04247087284264459

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf16"
$ws.Range("C3").Value = "Fgfr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1628603333333333
$ws.Range("H3").Value = 0.488581
$ws.Range("I3").Value = 0.06904471801498467
$ws.Range("J3").Value = 0.06904471801498467
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.5648773333333333
$ws.Range("N3").Value = 1.694632
$ws.Range("O3").Value = 0.2098038941373262
$ws.Range("P3").Value = 0.2098038941373262
$ws.Range("Q3").Value = 0.09199611079911112
$ws.Range("R3").Value = 0.8279649971919999
$ws.Range("S3").Value = 0.01448585070915738
$ws.Range("T3").Value = 0.01448585070915738

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf16"
$ws.Range("C4").Value = "Fgfr3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1628603333333333
$ws.Range("H4").Value = 0.488581
$ws.Range("I4").Value = 0.06904471801498467
$ws.Range("J4").Value = 0.06904471801498467
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4713726666666667
$ws.Range("N4").Value = 1.414118
$ws.Range("O4").Value = 0.1750748617810164
$ws.Range("P4").Value = 0.1750748617810165
$ws.Range("Q4").Value = 0.07676790961755556
$ws.Range("R4").Value = 0.6909111865580001
$ws.Range("S4").Value = 0.0120879944631827
$ws.Range("T4").Value = 0.0120879944631827

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fgf16"
$ws.Range("C5").Value = "Fgfr3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.769244333333333
$ws.Range("H5").Value = 5.307733
$ws.Range("I5").Value = 0.7500720009247772
$ws.Range("J5").Value = 0.7500720009247773
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.656156333333334
$ws.Range("N5").Value = 4.968469000000001
$ws.Range("O5").Value = 0.6151212440816572
$ws.Range("P5").Value = 0.6151212440816572
$ws.Range("Q5").Value = 2.930145207864111
$ws.Range("R5").Value = 26.371306870777
$ws.Range("S5").Value = 0.4613852223596669
$ws.Range("T5").Value = 0.461385222359667

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf16"
$ws.Range("C6").Value = "Fgfr3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.769244333333333
$ws.Range("H6").Value = 5.307733
$ws.Range("I6").Value = 0.7500720009247772
$ws.Range("J6").Value = 0.7500720009247773
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.5648773333333333
$ws.Range("N6").Value = 1.694632
$ws.Range("O6").Value = 0.2098038941373262
$ws.Range("P6").Value = 0.2098038941373262
$ws.Range("Q6").Value = 0.9994060210284443
$ws.Range("R6").Value = 8.994654189256
$ws.Range("S6").Value = 0.1573680266773944
$ws.Range("T6").Value = 0.1573680266773944

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf16"
$ws.Range("C7").Value = "Fgfr3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.769244333333333
$ws.Range("H7").Value = 5.307733
$ws.Range("I7").Value = 0.7500720009247772
$ws.Range("J7").Value = 0.7500720009247773
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.4713726666666667
$ws.Range("N7").Value = 1.414118
$ws.Range("O7").Value = 0.1750748617810164
$ws.Range("P7").Value = 0.1750748617810165
$ws.Range("Q7").Value = 0.8339734193882222
$ws.Range("R7").Value = 7.505760774494001
$ws.Range("S7").Value = 0.1313187518877158
$ws.Range("T7").Value = 0.1313187518877159

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Fgf16"
$ws.Range("C8").Value = "Fgfr3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.4266613333333333
$ws.Range("H8").Value = 1.279984
$ws.Range("I8").Value = 0.180883281060238
$ws.Range("J8").Value = 0.180883281060238
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.656156333333334
$ws.Range("N8").Value = 4.968469000000001
$ws.Range("O8").Value = 0.6151212440816572
$ws.Range("P8").Value = 0.6151212440816572
$ws.Range("Q8").Value = 0.7066178693884445
$ws.Range("R8").Value = 6.359560824496001
$ws.Range("S8").Value = 0.1112651488793457
$ws.Range("T8").Value = 0.1112651488793457

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Fgf16"
$ws.Range("C9").Value = "Fgfr3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.4266613333333333
$ws.Range("H9").Value = 1.279984
$ws.Range("I9").Value = 0.180883281060238
$ws.Range("J9").Value = 0.180883281060238
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.5648773333333333
$ws.Range("N9").Value = 1.694632
$ws.Range("O9").Value = 0.2098038941373262
$ws.Range("P9").Value = 0.2098038941373262
$ws.Range("Q9").Value = 0.2410113162097778
$ws.Range("R9").Value = 2.169101845888
$ws.Range("S9").Value = 0.03795001675077439
$ws.Range("T9").Value = 0.03795001675077439

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Fgf16"
$ws.Range("C10").Value = "Fgfr3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.4266613333333333
$ws.Range("H10").Value = 1.279984
$ws.Range("I10").Value = 0.180883281060238
$ws.Range("J10").Value = 0.180883281060238
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4713726666666667
$ws.Range("N10").Value = 1.414118
$ws.Range("O10").Value = 0.1750748617810164
$ws.Range("P10").Value = 0.1750748617810165
$ws.Range("Q10").Value = 0.2011164904568889
$ws.Range("R10").Value = 1.810048414112
$ws.Range("S10").Value = 0.03166811543011792
$ws.Range("T10").Value = 0.03166811543011792
